# Generate Report for Handoff
#
# The "Ready for handoff" rows (4-7) get refreshed when the handoff report
# regenerates:
#   - Priority goes from the placeholder "low" to the real "ht" value
#     (matching the already-handed-back rows above them), on the zh-cn
#     and de-de detail sheets.
#   - The Latest Handoff Datetime / Latest HO Xliff Generate Date
#     timestamps are bumped to the new generation time, on the detail
#     sheets as well as the Overview roll-up sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: rows 4-7, "Latest HO Xliff Generate Date" column
$overview.Range("G4:G7").Value = "2016-08-25 14:30:57"

# zh-cn sheet: rows 4-7
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4").Value = "2016-08-25 14:30:51"
$zhcn.Range("H5").Value = "2016-08-25 14:30:51"
$zhcn.Range("H6").Value = "2016-08-25 14:30:51"
$zhcn.Range("H7").Value = "2016-08-25 14:30:51"

# de-de sheet: rows 4-7
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4").Value = "2016-08-25 14:30:57"
$dede.Range("H5").Value = "2016-08-25 14:30:57"
$dede.Range("H6").Value = "2016-08-25 14:30:57"
$dede.Range("H7").Value = "2016-08-25 14:30:57"
